{"js": "// Replace the heading text \"Overcoming the Difficulties\" with \"Resources\",\n// preserving the run's existing formatting (bold, size, etc.).\nconst results = context.document.body.search(\"Overcoming the Difficulties\", {\n  matchCase: true,\n  matchWholeWord: false\n});\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Text \"Overcoming the Difficulties\" not found.');\n}\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"Resources\", Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Replace the heading text \"Overcoming the Difficulties\" with \"Resources\".\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"Overcoming the Difficulties\"\n$find.Replacement.Text = \"Resources\"\n$find.Forward = $true\n$find.Wrap = 1            # wdFindContinue\n$find.Format = $false\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n$find.MatchSoundsLike = $false\n$find.MatchAllWordForms = $false\n\n$find.Execute($find.Text, $false, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n"}
